$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark from its current location
#    (middle of paragraph 1, right after the " (woops! Lol)" run).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Add a brand-new paragraph right after paragraph 1.
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)

# 3. Type the first sentence into the new paragraph.
$r1 = $p2.Range
$r1.Collapse(1)
$r1.InsertAfter("No luck on round 1. Adding the diff to gitconfig…sooooo…maybe?")

# A temporary bookmark keeps the next chunk of text from being coalesced
# back into the previous run, so the saved document keeps the two runs
# separate (mirrors the two distinct <w:r> elements in the target).
$rMid = $d.Paragraphs(2).Range
$rMid.Collapse(0)
$d.Bookmarks.Add("zzzIronTempSplit", $rMid)

# 4. Type the second run of text.
$r2 = $d.Paragraphs(2).Range
$r2.Collapse(0)
$r2.InsertAfter(" NOOOOPE")

# Remove the temporary splitter bookmark now that the run break exists.
$d.Bookmarks("zzzIronTempSplit").Delete()

# 5. Re-add the "_GoBack" bookmark collapsed at the end of the new
#    paragraph's text (its last-edit-location semantics), matching the
#    diff. A collapsed range sitting immediately before a paragraph mark
#    confuses Bookmarks.Add, so temporarily append a placeholder
#    character, anchor the bookmark just before it, then remove the
#    placeholder again.
$rPad = $d.Paragraphs(2).Range
$rPad.Collapse(0)
$rPad.InsertAfter("X")

$padEnd = $d.Paragraphs(2).Range.End - 2
$rBm = $d.Range($padEnd, $padEnd)
$d.Bookmarks.Add("_GoBack", $rBm)

$delStart = $d.Paragraphs(2).Range.End - 2
$d.Range($delStart, $delStart + 1).Delete()

Write-Output "done"
